$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "29.406.00"
Set-TextCell $ws.Range("E2") "  +0.96%  "

Set-TextCell $ws.Range("D3") "1.948.90"
Set-TextCell $ws.Range("E3") "  +2.52%  "

Set-TextCell $ws.Range("D4") "1.005"
Set-TextCell $ws.Range("E4") "  +0.46%  "

Set-TextCell $ws.Range("D5") "326.39"
Set-TextCell $ws.Range("E5") "  +0.17%  "

Set-TextCell $ws.Range("D6") "1.004"
Set-TextCell $ws.Range("E6") "  +0.33%  "

Set-TextCell $ws.Range("D7") "0.4632"
Set-TextCell $ws.Range("E7") "  +0.16%  "

Set-TextCell $ws.Range("D8") "0.3874"
Set-TextCell $ws.Range("E8") "  -0.70%  "

Set-TextCell $ws.Range("D9") "46.27"
Set-TextCell $ws.Range("E9") "  +0.77%  "

Set-TextCell $ws.Range("D10") "0.07840"
Set-TextCell $ws.Range("E10") "  -0.51%  "

Set-TextCell $ws.Range("D11") "0.9809"
Set-TextCell $ws.Range("E11") "  -1.10%  "

Set-TextCell $ws.Range("D12") "22.74"
Set-TextCell $ws.Range("E12") "  +3.85%  "

Set-TextCell $ws.Range("D13") "1.946.87"
Set-TextCell $ws.Range("E13") "  +2.56%  "

Set-TextCell $ws.Range("D14") "7.104"
Set-TextCell $ws.Range("E14") "  +0.52%  "

Set-TextCell $ws.Range("D15") "5.756"

Set-TextCell $ws.Range("D16") "0.07068"
Set-TextCell $ws.Range("E16") "  +1.14%  "

Set-TextCell $ws.Range("D17") "87.05"
Set-TextCell $ws.Range("E17") "  -1.04%  "

Set-TextCell $ws.Range("D18") "1.007"
Set-TextCell $ws.Range("E18") "  +0.51%  "

Set-TextCell $ws.Range("D19") "0.000009859"
Set-TextCell $ws.Range("E19") "  -0.60%  "

Set-TextCell $ws.Range("D20") "17.01"
Set-TextCell $ws.Range("E20") "  -0.34%  "

Set-TextCell $ws.Range("D21") "1.004"
Set-TextCell $ws.Range("E21") "  +0.29%  "

Set-TextCell $ws.Range("D22") "29.417.27"
Set-TextCell $ws.Range("E22") "  +0.94%  "

Set-TextCell $ws.Range("D23") "5.488"
Set-TextCell $ws.Range("E23") "  +3.14%  "

Set-TextCell $ws.Range("D24") "11.08"
Set-TextCell $ws.Range("E24") "  -0.33%  "

Set-TextCell $ws.Range("D25") "2.185.18"
Set-TextCell $ws.Range("E25") "  +2.89%  "

Set-TextCell $ws.Range("D26") "2.103"
Set-TextCell $ws.Range("E26") "  -0.35%  "

Set-TextCell $ws.Range("D27") "157.56"
Set-TextCell $ws.Range("E27") "  +1.07%  "

Set-TextCell $ws.Range("E28") "  -0.03%  "

Set-TextCell $ws.Range("D29") "5.775"
Set-TextCell $ws.Range("E29") "  -2.43%  "

Set-TextCell $ws.Range("D30") "118.63"
Set-TextCell $ws.Range("E30") "  -0.06%  "

Set-TextCell $ws.Range("D31") "1.861"
Set-TextCell $ws.Range("E31") "  -0.99%  "

Set-TextCell $ws.Range("D32") "0.09390"
Set-TextCell $ws.Range("E32") "  +0.64%  "

Set-TextCell $ws.Range("D33") "0.8644"
Set-TextCell $ws.Range("E33") "  -3.94%  "

Set-TextCell $ws.Range("D34") "5.179"
Set-TextCell $ws.Range("E34") "  -1.44%  "

Set-TextCell $ws.Range("D35") "1.302"
Set-TextCell $ws.Range("E35") "  -1.75%  "

Set-TextCell $ws.Range("D36") "3.128"
Set-TextCell $ws.Range("E36") "  -0.95%  "

Set-TextCell $ws.Range("D37") "0.05755"
Set-TextCell $ws.Range("E37") "  -0.88%  "

Set-TextCell $ws.Range("B38") "VeChain"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D38") "0.02092"
Set-TextCell $ws.Range("E38") "  +0.26%  "

Set-TextCell $ws.Range("B39") "TrustWalletToken"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Range("D39") "1.148"
Set-TextCell $ws.Range("E39") "  -2.29%  "

Set-TextCell $ws.Range("D40") "7.700"
Set-TextCell $ws.Range("E40") "  -0.38%  "

Set-TextCell $ws.Range("D41") "0.5670"
Set-TextCell $ws.Range("E41") "  -0.42%  "

Set-TextCell $ws.Range("D42") "0.1783"
Set-TextCell $ws.Range("E42") "  -0.56%  "

Set-TextCell $ws.Range("D43") "9.437"
Set-TextCell $ws.Range("E43") "  -3.17%  "

Set-TextCell $ws.Range("B44") "PEPE"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D44") "0.000002860"
Set-TextCell $ws.Range("E44") "  +49.02%  "

Set-TextCell $ws.Range("B45") "MXToken"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws.Range("D45") "2.736"
Set-TextCell $ws.Range("E45") "  +7.19%  "

Set-TextCell $ws.Range("B46") "Decentraland"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws.Range("D46") "0.5298"
Set-TextCell $ws.Range("E46") "  -1.18%  "

Set-TextCell $ws.Range("B47") "EnergySwap"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D47") "11.61"
Set-TextCell $ws.Range("E47") "  -2.43%  "

Set-TextCell $ws.Range("D48") "2.132"
Set-TextCell $ws.Range("E48") "  -4.94%  "

Set-TextCell $ws.Range("D49") "0.06875"
Set-TextCell $ws.Range("E49") "  -1.92%  "

Set-TextCell $ws.Range("D50") "1.818"
Set-TextCell $ws.Range("E50") "  -1.68%  "

Set-TextCell $ws.Range("D51") "111.82"
Set-TextCell $ws.Range("E51") "  -1.18%  "
